# Add I0 and IF columns (I and J) to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF", styled like the other header cells (copy style from H1)
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Data rows 2-92 for columns I (9) and J (10)
$data = @(
    @{Row=2; I=2; J=2},
    @{Row=3; I=6; J=6},
    @{Row=4; I=11; J=11},
    @{Row=5; I=7; J=7},
    @{Row=6; I=4; J=5},
    @{Row=7; I=7; J=7},
    @{Row=8; I=8; J=8},
    @{Row=9; I=4; J=4},
    @{Row=10; I=6; J=6},
    @{Row=11; I=8; J=8},
    @{Row=12; I=4; J=5},
    @{Row=13; I=8; J=8},
    @{Row=14; I=6; J=6},
    @{Row=15; I=6; J=6},
    @{Row=16; I=6; J=6},
    @{Row=17; I=7; J=7},
    @{Row=18; I=6; J=7},
    @{Row=19; I=8; J=8},
    @{Row=20; I=9; J=9},
    @{Row=21; I=6; J=6},
    @{Row=22; I=5; J=5},
    @{Row=23; I=7; J=7},
    @{Row=24; I=10; J=10},
    @{Row=25; I=8; J=8},
    @{Row=26; I=7; J=7},
    @{Row=27; I=9; J=9},
    @{Row=28; I=8; J=8},
    @{Row=29; I=6; J=6},
    @{Row=30; I=7; J=7},
    @{Row=31; I=6; J=6},
    @{Row=32; I=7; J=7},
    @{Row=33; I=8; J=8},
    @{Row=34; I=9; J=9},
    @{Row=35; I=8; J=8},
    @{Row=36; I=7; J=7},
    @{Row=37; I=9; J=9},
    @{Row=38; I=7; J=7},
    @{Row=39; I=8; J=8},
    @{Row=40; I=8; J=8},
    @{Row=41; I=8; J=8},
    @{Row=42; I=8; J=8},
    @{Row=43; I=7; J=7},
    @{Row=44; I=6; J=6},
    @{Row=45; I=7; J=7},
    @{Row=46; I=9; J=9},
    @{Row=47; I=7; J=7},
    @{Row=48; I=7; J=7},
    @{Row=49; I=8; J=8},
    @{Row=50; I=7; J=7},
    @{Row=51; I=7; J=7},
    @{Row=52; I=8; J=8},
    @{Row=53; I=8; J=8},
    @{Row=54; I=8; J=9},
    @{Row=55; I=9; J=9},
    @{Row=56; I=9; J=9},
    @{Row=57; I=9; J=9},
    @{Row=58; I=9; J=9},
    @{Row=59; I=9; J=9},
    @{Row=60; I=9; J=9},
    @{Row=61; I=9; J=9},
    @{Row=62; I=7; J=7},
    @{Row=63; I=9; J=9},
    @{Row=64; I=9; J=9},
    @{Row=65; I=8; J=8},
    @{Row=66; I=9; J=9},
    @{Row=67; I=9; J=9},
    @{Row=68; I=9; J=9},
    @{Row=69; I=9; J=9},
    @{Row=70; I=9; J=9},
    @{Row=71; I=8; J=8},
    @{Row=72; I=8; J=8},
    @{Row=73; I=8; J=8},
    @{Row=74; I=8; J=8},
    @{Row=75; I=7; J=8},
    @{Row=76; I=8; J=8},
    @{Row=77; I=9; J=9},
    @{Row=78; I=7; J=7},
    @{Row=79; I=8; J=8},
    @{Row=80; I=8; J=8},
    @{Row=81; I=7; J=7},
    @{Row=82; I=8; J=9},
    @{Row=83; I=7; J=7},
    @{Row=84; I=7; J=8},
    @{Row=85; I=8; J=8},
    @{Row=86; I=9; J=9},
    @{Row=87; I=8; J=8},
    @{Row=88; I=9; J=9},
    @{Row=89; I=5; J=5},
    @{Row=90; I=6; J=6},
    @{Row=91; I=3; J=3},
    @{Row=92; I=5; J=5}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
}
